$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells (revised AgTests/AgPosit figures) ---
$ws.Range("F473").Value = 40050
$ws.Range("F479").Value = 42633
$ws.Range("F480").Value = 33583
$ws.Range("F483").Value = 66162
$ws.Range("F486").Value = 8977
$ws.Range("F487").Value = 6902
$ws.Range("F488").Value = 6431
$ws.Range("F493").Value = 8376
$ws.Range("F494").Value = 6800
$ws.Range("F498").Value = 9276
$ws.Range("F500").Value = 7864
$ws.Range("F501").Value = 5863
$ws.Range("F507").Value = 7439
$ws.Range("F508").Value = 5872
$ws.Range("F514").Value = 7195
$ws.Range("F515").Value = 5219
$ws.Range("F521").Value = 6972
$ws.Range("F522").Value = 5250
$ws.Range("F528").Value = 8180
$ws.Range("F529").Value = 5808
$ws.Range("F530").Value = 12955
$ws.Range("F531").Value = 9337
$ws.Range("F533").Value = 11917
$ws.Range("F535").Value = 10199
$ws.Range("F536").Value = 8014
$ws.Range("G538").Value = 29
$ws.Range("F539").Value = 10648
$ws.Range("F542").Value = 10369
$ws.Range("F543").Value = 4721
$ws.Range("F544").Value = 14378
$ws.Range("F548").Value = 17169
$ws.Range("F550").Value = 8505
$ws.Range("F551").Value = 17822
$ws.Range("F552").Value = 15594
$ws.Range("F553").Value = 15413
$ws.Range("F554").Value = 17908
$ws.Range("F558").Value = 24937
$ws.Range("F559").Value = 22483
$ws.Range("F561").Value = 24045
$ws.Range("G561").Value = 393
$ws.Range("F562").Value = 27015
$ws.Range("G562").Value = 278
$ws.Range("F563").Value = 14066
$ws.Range("G563").Value = 174
$ws.Range("F564").Value = 14247
$ws.Range("G564").Value = 197
$ws.Range("F565").Value = 28742
$ws.Range("G565").Value = 368
$ws.Range("F566").Value = 25633
$ws.Range("G566").Value = 326
$ws.Range("F567").Value = 23267
$ws.Range("G567").Value = 319
$ws.Range("F568").Value = 23348
$ws.Range("G568").Value = 290
$ws.Range("F569").Value = 31781
$ws.Range("G569").Value = 353
$ws.Range("F570").Value = 14746
$ws.Range("G570").Value = 219
$ws.Range("F571").Value = 15253
$ws.Range("G571").Value = 277

# --- Append new daily rows ---
$ws.Range("A572").Value = 44466
$ws.Range("B572").Value = 409621
$ws.Range("C572").Value = 11086
$ws.Range("D572").Value = 1012
$ws.Range("E572").Value = 12606
$ws.Range("F572").Value = 32123
$ws.Range("G572").Value = 576
$ws.Range("A573").Value = 44467
$ws.Range("B573").Value = 411080
$ws.Range("C573").Value = 11698
$ws.Range("D573").Value = 1459
$ws.Range("E573").Value = 12620
$ws.Range("F573").Value = 25155
$ws.Range("G573").Value = 376
$ws.Range("A574").Value = 44468
$ws.Range("B574").Value = 412507
$ws.Range("C574").Value = 11174
$ws.Range("D574").Value = 1427
$ws.Range("E574").Value = 12637
$ws.Range("F574").Value = 15547
$ws.Range("G574").Value = 225


